# Add a new method to TestData class.
# - Set TestInputs (2nd sheet) as the active/selected sheet
# - Update A7 label from "Order" to "Order ID"
# - Move the selection on TestInputs to B7

$wb = $excel.ActiveWorkbook

$wsTestInputs = $wb.Worksheets.Item("TestInputs")

# Update the label in A7
$wsTestInputs.Range("A7").Value = "Order ID"

# Activate the TestInputs sheet and select B7, making it the active tab/selection
$wsTestInputs.Activate()
$wsTestInputs.Range("B7").Select()
